$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.447.56'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '3.702.93'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '613.18'
$ws.Range('E5').Value = '  +6.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '195.62'
$ws.Range('E6').Value = '  +14.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.634'
$ws.Range('E7').Value = '  +2.32%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.724'
$ws.Range('E9').Value = '  +3.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '60.81'
$ws.Range('E10').Value = '  +19.28%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000286'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.43'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '4.294.01'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = '3.704.68'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.47'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.127'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.15'
$ws.Range('E18').Value = '  +3.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.84'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('D20').Value = '68.326.25'
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '408.56'
$ws.Range('E21').Value = '  +1.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.62'
$ws.Range('E22').Value = '  +3.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '89.91'
$ws.Range('E23').Value = '  +3.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.54'
$ws.Range('E24').Value = '  +9.16%  '
$ws.Range('E25').Value = '  +2.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.11'
$ws.Range('E26').Value = '  +3.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.03'
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.77'
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  +3.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.74'
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.72'
$ws.Range('E31').Value = '  +4.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '48.23'
$ws.Range('E32').Value = '  +12.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.69'
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.122'
$ws.Range('E34').Value = '  +5.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '635.72'
$ws.Range('E35').Value = '  +7.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '67.38'
$ws.Range('E36').Value = '  +3.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.413'
$ws.Range('E37').Value = '  +5.16%  '
$ws.Range('D38').Value = '0.0₃0811'
$ws.Range('E38').Value = '  -7.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('E41').Value = '  +5.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.04'
$ws.Range('E42').Value = '  +2.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0444'
$ws.Range('E43').Value = '  +2.68%  '
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '2.929.10'
$ws.Range('E45').Value = '  +5.22%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.140'
$ws.Range('E46').Value = '  +5.41%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.37'
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.69'
$ws.Range('E48').Value = '  +1.63%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '145.89'
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.66'
$ws.Range('E50').Value = '  -7.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.05'
$ws.Range('E51').Value = '  -3.28%  '
